$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-01-06 Monday" "2025-01-07 Tuesday"
Replace-Text "879×9=7911" "146×9=1314"
Replace-Text "204×3=612" "295×7=2065"
Replace-Text "821×6=4926" "351×6=2106"
Replace-Text "336×4=1344" "544×6=3264"
Replace-Text "569×7=3983" "160×8=1280"
Replace-Text "772×7=5404" "154×2=308"
Replace-Text "335×4=1340" "441×6=2646"
Replace-Text "673×3=2019" "443×6=2658"
Replace-Text "301×4=1204" "532×4=2128"
Replace-Text "830×5=4150" "540×8=4320"
Replace-Text "953×9=8577" "773×6=4638"
Replace-Text "804×6=4824" "585×6=3510"
Replace-Text "793×2=1586" "139×3=417"
Replace-Text "640×6=3840" "761×8=6088"
Replace-Text "207×9=1863" "899×7=6293"
Replace-Text "759×6=4554" "163×5=815"
Replace-Text "363×4=1452" "144×4=576"
Replace-Text "657×9=5913" "509×6=3054"
Replace-Text "323×3=969" "108×7=756"
Replace-Text "404×2=808" "511×9=4599"
Replace-Text "500×6=3000" "472×9=4248"
Replace-Text "625×8=5000" "317×4=1268"
Replace-Text "285×3=855" "881×3=2643"
Replace-Text "121×3=363" "174×3=522"
Replace-Text "192×8=1536" "678×6=4068"
